# Insert a new leading column (A) on sheet1.
#
# Before: col A held the question text, col B held the answer text.
# After:  col A holds a 0-based row index number, col B holds the
#         (shifted) question text, and col C holds the (shifted) answer
#         text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Shift the existing two columns of data one column to the right,
# leaving a blank column A in their place.
$ws.Columns.Item(1).Insert()

# Column B (the original column A) now holds the question text for every
# data row; walk down it to find how many rows of data exist.
$lastRow = $ws.Cells.Item(1, 2).End(4).Row

# Fill the new column A with sequential numbers, starting at 0.
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
